$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.464.82'
$ws.Range('E2').Value = '  +0.13%  '

$ws.Range('D3').Value = '3.524.48'
$ws.Range('E3').Value = '  +0.08%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '615.62'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.37%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '151.63'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.44%  '

$ws.Range('D7').Value = '3.523.91'
$ws.Range('E7').Value = '  +0.18%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.482'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.52%  '

$ws.Range('E10').Value = '  -0.71%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.13'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.55%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.425'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.46%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000221'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.31%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '32.10'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.55%  '

$ws.Range('D15').Value = '4.121.00'
$ws.Range('E15').Value = '  +0.15%  '

$ws.Range('D16').Value = '3.529.67'
$ws.Range('E16').Value = '  +0.73%  '

$ws.Range('D17').Value = '67.473.55'
$ws.Range('E17').Value = '  +0.10%  '

$ws.Range('E18').Value = '  -0.15%  '

$ws.Range('E19').Value = '  +0.78%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '15.40'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.24%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '446.37'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.73%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '9.54'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.68%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.624'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.20%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '77.55'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.34%  '

$ws.Range('E25').Value = '  +11.24%  '

$ws.Range('D26').Value = '3.666.09'
$ws.Range('E26').Value = '  +0.16%  '

$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.24'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.61%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.48'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.59%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.52'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.10%  '

$ws.Range('E31').Value = '  -8.44%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.04%  '

$ws.Range('E33').Value = '  +4.24%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '25.86'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.58%  '

$ws.Range('E35').Value = '  -0.83%  '

$ws.Range('D36').Value = '3.517.29'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('E37').Value = '  -3.48%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.02'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.00%  '

$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.01%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '177.41'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.44%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.16'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +3.26%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0883'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +0.34%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '5.44'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.36%  '

$ws.Range('E45').Value = '  -0.75%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '28.40'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.81%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '45.12'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.28%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.65'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.71%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.25'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +3.31%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.59'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.88%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.999'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.73%  '
